$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = 0.17344416360779746
$ws.Range("A2").Value = -0.00599999999543499
$ws.Range("A3").Value = -0.003999999996063153
$ws.Range("A4").Value = -0.00799999999280665
$ws.Range("A5").Value = -0.002999999996183611
$ws.Range("A6").Value = -0.001999999996257884
$ws.Range("A7").Value = -0.00999999999031509
$ws.Range("A8").Value = -0.009999999990179642
$ws.Range("A9").Value = 0.022381441216825504
$ws.Range("A10").Value = -0.00199999999568945
$ws.Range("A11").Value = -0.002999999994888647
$ws.Range("A12").Value = -0.0034999999944154148
$ws.Range("A13").Value = -0.0034999999940659166
$ws.Range("A14").Value = -0.007999999990667028
$ws.Range("A15").Value = -0.0009999999957681638
$ws.Range("A16").Value = -0.0019999999949669167
$ws.Range("A17").Value = -0.0019999999948820957
$ws.Range("A18").Value = -0.003999999993389736
$ws.Range("A19").Value = -0.003999999996891823
$ws.Range("A20").Value = -0.003999999996693759
$ws.Range("A21").Value = -0.003999999996636028
$ws.Range("A22").Value = -0.0039999999966244815
$ws.Range("A23").Value = -0.004999999995046522
$ws.Range("A24").Value = -0.034453546048953676
$ws.Range("A25").Value = -0.019999999983227212
$ws.Range("A26").Value = -0.002499999995745128
$ws.Range("A27").Value = -0.0024999999956474284
$ws.Range("A28").Value = -0.016857821416547125
$ws.Range("A29").Value = -0.006999999991652572
$ws.Range("A30").Value = -0.05999999995258731
$ws.Range("A31").Value = -0.00699999999210732
$ws.Range("A32").Value = -0.0064034339292415865
$ws.Range("A33").Value = 0.060024928256360965
